$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Prada ---
$ws.Range("C4").Value = 55.35
$ws.Range("D4").Formula = "=+C4*0.85"

# --- Row 5: Apple ---
$ws.Range("C5").Value = 234.96
$ws.Range("D5").Value = 200

# --- Row 6: Adobe ---
$ws.Range("C6").Value = 496.29
$ws.Range("D6").Value = 307.6998

# --- Row 13: ASML ---
$ws.Range("C13").Value = 720
$ws.Range("D13").Value = 200

# --- Row 59 (new row): Ferroglobe - create percent style (must be created before the italic $ style below) ---
$ws.Range("D59").Value = $null
$ws.Range("D59").NumberFormat = "0.00%"

# --- Row 45: RocketLab - empty cell, italic $ style ---
$ws.Range("D45").Value = $null
$ws.Range("D45").Font.Italic = $true

# --- Row 50: Toyota Motor Corporation ---
$ws.Range("C50").Value = 171.72

# --- Row 57: Baxter International (set before row 54 so "<35" gets the lower shared-string index) ---
$ws.Range("C57").Value = 36.55
$ws.Range("D57").Value = "<35"

# --- Row 54: Zegna ---
$ws.Range("C54").Value = 8.4
$ws.Range("D54").Value = "<8"

# --- Row 55: T-Mobile - replace formula with static value ---
$ws.Range("D55").Value = 150

# --- Row 58: ICU Medical ---
$ws.Range("C58").Value = 183.05
$ws.Range("D58").Value = 54.863725490196096
$ws.Range("D58").Borders.LineStyle = -4142

# --- Row 59 (new row): Ferroglobe - set the name now ---
$ws.Range("B59").Value = "Ferroglobe"

# --- Update selection to match saved workbook state ---
$ws.Range("C62").Select()
